$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.643.30"
$ws.Range("E2").Value = "  -1.25%  "
$ws.Range("D3").Value = "1.592.61"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.97"
$ws.Range("E5").Value = "  -1.81%  "
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  -1.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0617"
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.67"
$ws.Range("E10").Value = "  -2.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0834"
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("D12").Value = "1.813.75"
$ws.Range("E12").Value = "  -1.83%  "
$ws.Range("D13").Value = "1.600.40"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.04"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("E15").Value = "  -2.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.71"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "26.629.47"
$ws.Range("E17").Value = "  -1.36%  "
$ws.Range("E18").Value = "  -1.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.87"
$ws.Range("E19").Value = "  -3.16%  "
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.74"
$ws.Range("E21").Value = "  -2.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.26"
$ws.Range("E22").Value = "  -2.02%  "
$ws.Range("E23").Value = "  -2.31%  "
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.83"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("E27").Value = "  -0.06%  "
$ws.Range("E28").Value = "  -2.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.32"
$ws.Range("E29").Value = "  -1.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0508"
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("E31").Value = "  -1.77%  "
$ws.Range("E32").Value = "  -3.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.665"
$ws.Range("E33").Value = "  +22.54%  "
$ws.Range("E34").Value = "  -2.23%  "
$ws.Range("D35").Value = "1.320.76"
$ws.Range("E35").Value = "  -0.94%  "
$ws.Range("E36").Value = "  -4.01%  "
$ws.Range("E37").Value = "  -2.26%  "
$ws.Range("E38").Value = "  -1.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.832"
$ws.Range("E39").Value = "  -1.51%  "
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("E41").Value = "  +3.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.788"
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.16"
$ws.Range("E44").Value = "  -1.97%  "
$ws.Range("D45").Value = "1.726.67"
$ws.Range("E45").Value = "  -1.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.04"
$ws.Range("E46").Value = "  -0.39%  "
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.839"
$ws.Range("E48").Value = "  +2.33%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0980"
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.51"
$ws.Range("E51").Value = "  -0.30%  "
